$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge "THU Oct 12" / " 11:10:37 PDT 2017" into a single run/text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("THU Oct 12 11:10:37 PDT 2017", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "THU Oct 12 11:10:37 PDT 2017", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Append the new "FRI Oct 13" purchase-details block after the
#    "Amount balance ... - 27282.0" paragraph.
# ---------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("- 27282.0", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $findRange.Paragraphs(1)
$anchorIndex = $anchorPara.Range.Information(3)  # wdActiveEndAdjustedPageNumber placeholder (unused)

# Locate the numeric paragraph index so we can come back and format things.
$startCount = $d.Paragraphs.Count
$insertPoint = $d.Range($anchorPara.Range.End, $anchorPara.Range.End)

$lines = @(
    "",
    "FRI Oct 13 11:16:57 PDT 2017",
    ("Person Name" + ([char]9) + ([char]9) + ([char]9) + ([char]9) + "- NS"),
    "---------------------------------------------------------------",
    ("Item Name" + ([char]9) + ([char]9) + ([char]9) + ([char]9) + "- CARROT EVE"),
    ("Number of Pockets" + ([char]9) + ([char]9) + ([char]9) + "- 1"),
    ("Number of KGs" + ([char]9) + ([char]9) + ([char]9) + "- 76"),
    ("Rate" + ([char]9) + ([char]9) + ([char]9) + ([char]9) + ([char]9) + "- 22"),
    ("Total Price" + ([char]9) + ([char]9) + ([char]9) + ([char]9) + "- 1672.0"),
    ("Amount balance" + ([char]9) + ([char]9) + ([char]9) + "- 28954.0"),
    "",
    "",
    ""
)

$blockText = [string]::Join("`r", $lines)
$insertPoint.InsertAfter($blockText)

# Paragraph index (1-based) of the first newly inserted paragraph (the
# blank bold one right after "...- 27282.0").
$firstNewIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $anchorPara.Range.End) {
        $firstNewIndex = $i
        break
    }
}

# Bold paragraphs, relative to $firstNewIndex (0-based offsets):
#   0  -> blank bold paragraph
#   9  -> "Amount balance ... - 28954.0"
#   10 -> blank bold paragraph
#   11 -> blank bold paragraph
#   12 -> blank bold paragraph
$boldOffsets = @(0, 9, 10, 11, 12)
foreach ($off in $boldOffsets) {
    $p = $d.Paragraphs($firstNewIndex + $off)
    $p.Range.Font.Bold = 1
}

# Make sure all the other (non-bold) new paragraphs are explicitly not bold.
for ($off = 1; $off -le 8; $off++) {
    $p = $d.Paragraphs($firstNewIndex + $off)
    $p.Range.Font.Bold = 0
}

# Ensure every newly added paragraph uses the "Plain Text" style, matching
# the rest of the document.
for ($off = 0; $off -le 12; $off++) {
    $p = $d.Paragraphs($firstNewIndex + $off)
    $p.Style = $d.Styles("Plain Text")
}

Write-Host "Done. Paragraphs before:" $startCount "after:" $d.Paragraphs.Count
